# Rotate the contents of columns B, C and D one column to the left
# (new B = old C, new C = old D, new D = old B) for the populated rows
# of Sheet1, then swap the bestFit column widths so the formatting
# keeps following the data, and finally restore the saved selection
# state (the whole of column C).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellValuePreserveFormat($cell, $value) {
    # Cells in this sheet are formatted as Text ("@"). Writing a numeric
    # value into a Text-formatted cell gets stored as text instead of a
    # real number, which would diverge from the source data. Flip the
    # destination to General just long enough to write a true numeric
    # value, then restore its original number format.
    $isNumeric = $value -is [double] -or $value -is [int]
    if ($isNumeric -and $cell.NumberFormat -eq "@") {
        $savedFormat = $cell.NumberFormat
        $cell.NumberFormat = "General"
        $cell.Value2 = $value
        $cell.NumberFormat = $savedFormat
    } else {
        $cell.Value2 = $value
    }
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($row = 1; $row -le $lastRow; $row++) {
    $bCell = $ws.Cells.Item($row, 2)
    $cCell = $ws.Cells.Item($row, 3)
    $dCell = $ws.Cells.Item($row, 4)

    $bVal = $bCell.Value2
    $cVal = $cCell.Value2
    $dVal = $dCell.Value2

    Set-CellValuePreserveFormat $bCell $cVal
    Set-CellValuePreserveFormat $cCell $dVal
    Set-CellValuePreserveFormat $dCell $bVal
}

# Column B/C now hold the (narrower) data that used to live in C/D, and
# column D holds the (wider) data that used to live in B, so swap the
# bestFit column widths to match (ColumnWidth is expressed in characters
# and gets snapped to whole pixels, so these are the values that land on
# the same pixel widths as the original B/(C:D) columns).
$ws.Columns.Item(2).ColumnWidth = 6
$ws.Columns.Item(3).ColumnWidth = 6
$ws.Columns.Item(4).ColumnWidth = 4.14

# Restore the selection that was captured on save: the whole of column C.
$ws.Range("C1:C1048576").Select() | Out-Null
